# Updates cryptos list price/volume(1h) columns (D, E) for Sat Jan 20 2024 run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextSafeValue($range, [string]$text) {
    # Excel auto-converts plain-decimal-looking strings (e.g. "0.999") to numbers,
    # which would lose the original text formatting / trailing zeros (e.g. "9.70" -> 9.7).
    # Force those through as text with a quote-prefix, exactly as typing '0.999 would in Excel,
    # while leaving already-unambiguous text (dates-like "41.526.50", unicode subscripts, etc.) alone.
    if ($text -match '^\d+\.\d+$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

Set-TextSafeValue $ws.Range("D2") "41.526.50"
$ws.Range("E2").Value = "  +0.03%  "
Set-TextSafeValue $ws.Range("D3") "2.470.42"
$ws.Range("E3").Value = "  -0.68%  "
Set-TextSafeValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.24%  "
Set-TextSafeValue $ws.Range("D5") "314.89"
$ws.Range("E5").Value = "  +0.39%  "
Set-TextSafeValue $ws.Range("D6") "92.09"
$ws.Range("E6").Value = "  -2.65%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +2.87%  "
Set-TextSafeValue $ws.Range("D10") "32.35"
$ws.Range("E10").Value = "  -3.81%  "
Set-TextSafeValue $ws.Range("D11") "0.0792"
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("E12").Value = "  +0.91%  "
Set-TextSafeValue $ws.Range("D13") "2.847.62"
$ws.Range("E13").Value = "  -0.73%  "
Set-TextSafeValue $ws.Range("D14") "6.85"
$ws.Range("E14").Value = "  -2.07%  "
Set-TextSafeValue $ws.Range("D15") "15.98"
$ws.Range("E15").Value = "  +3.03%  "
Set-TextSafeValue $ws.Range("D16") "2.503.99"
$ws.Range("E16").Value = "  +2.59%  "
Set-TextSafeValue $ws.Range("D17") "0.776"
$ws.Range("E17").Value = "  -2.24%  "
Set-TextSafeValue $ws.Range("D18") "41.543.32"
$ws.Range("E18").Value = "  +0.17%  "
Set-TextSafeValue $ws.Range("D19") "6.49"
$ws.Range("E19").Value = "  +1.99%  "
Set-TextSafeValue $ws.Range("D20") "0.0₃0943"
$ws.Range("E20").Value = "  +2.01%  "
Set-TextSafeValue $ws.Range("D21") "71.12"
$ws.Range("E21").Value = "  +3.33%  "
Set-TextSafeValue $ws.Range("D22") "11.11"
$ws.Range("E22").Value = "  -1.80%  "
Set-TextSafeValue $ws.Range("D23") "236.84"
$ws.Range("E23").Value = "  -0.16%  "
Set-TextSafeValue $ws.Range("D24") "2.72"
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("E26").Value = "  -0.51%  "
Set-TextSafeValue $ws.Range("D27") "24.74"
$ws.Range("E27").Value = "  +2.31%  "
Set-TextSafeValue $ws.Range("D28") "2.24"
$ws.Range("E28").Value = "  -0.14%  "
Set-TextSafeValue $ws.Range("D29") "9.70"
$ws.Range("E29").Value = "  -0.77%  "
Set-TextSafeValue $ws.Range("D30") "35.43"
$ws.Range("E30").Value = "  -3.46%  "
Set-TextSafeValue $ws.Range("D31") "155.92"
$ws.Range("E31").Value = "  +2.40%  "
Set-TextSafeValue $ws.Range("D32") "5.45"
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("E34").Value = "  +0.11%  "
Set-TextSafeValue $ws.Range("D35") "17.27"
$ws.Range("E35").Value = "  -4.89%  "
Set-TextSafeValue $ws.Range("D36") "2.88"
$ws.Range("E36").Value = "  -6.94%  "
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("E38").Value = "  +1.23%  "
Set-TextSafeValue $ws.Range("D39") "1.78"
$ws.Range("E39").Value = "  -5.41%  "
Set-TextSafeValue $ws.Range("D40") "2.22"
$ws.Range("E40").Value = "  -11.14%  "
Set-TextSafeValue $ws.Range("D41") "4.02"
$ws.Range("E41").Value = "  -5.36%  "
$ws.Range("E42").Value = "  -0.31%  "
Set-TextSafeValue $ws.Range("D43") "1.946.83"
$ws.Range("E43").Value = "  -2.27%  "
$ws.Range("E44").Value = "  -1.26%  "
Set-TextSafeValue $ws.Range("D45") "18.78"
$ws.Range("E45").Value = "  -5.19%  "
Set-TextSafeValue $ws.Range("D46") "2.94"
$ws.Range("E46").Value = "  -3.13%  "
Set-TextSafeValue $ws.Range("D47") "9.05"
$ws.Range("E47").Value = "  +2.11%  "
Set-TextSafeValue $ws.Range("D48") "2.706.35"
$ws.Range("E48").Value = "  -0.97%  "
Set-TextSafeValue $ws.Range("D49") "96.97"
$ws.Range("E49").Value = "  -0.18%  "
Set-TextSafeValue $ws.Range("D50") "67.31"
$ws.Range("E50").Value = "  -3.77%  "
Set-TextSafeValue $ws.Range("D51") "52.38"
$ws.Range("E51").Value = "  +2.64%  "
